# Apply updated TPM-derived values to LR-pairs sheet (F12-Gp1ba)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.7822470115268996
$ws.Range("J2").Value = 0.7822470115268996
$ws.Range("M2").Value = 3.704480666666667
$ws.Range("N2").Value = 11.113442
$ws.Range("O2").Value = 0.298964201395561
$ws.Range("P2").Value = 0.2989642013955609
$ws.Range("Q2").Value = 0.5067309710857778
$ws.Range("R2").Value = 4.560578739772
$ws.Range("S2").Value = 0.2338638530952037
$ws.Range("T2").Value = 0.2338638530952037

# Row 3
$ws.Range("I3").Value = 0.7822470115268996
$ws.Range("J3").Value = 0.7822470115268996
$ws.Range("O3").Value = 0.3765624616238499
$ws.Range("P3").Value = 0.3765624616238499
$ws.Range("S3").Value = 0.2945648602584694
$ws.Range("T3").Value = 0.2945648602584694

# Row 4
$ws.Range("I4").Value = 0.7822470115268996
$ws.Range("J4").Value = 0.7822470115268996
$ws.Range("M4").Value = 2.870093333333334
$ws.Range("N4").Value = 8.610280000000001
$ws.Range("O4").Value = 0.2316263029934534
$ws.Range("P4").Value = 0.2316263029934534
$ws.Range("Q4").Value = 0.3925962402755556
$ws.Range("R4").Value = 3.533366162480001
$ws.Range("S4").Value = 0.1811889833076531
$ws.Range("T4").Value = 0.1811889833076531

# Row 5
$ws.Range("I5").Value = 0.7822470115268996
$ws.Range("J5").Value = 0.7822470115268996
$ws.Range("M5").Value = 1.150472333333333
$ws.Range("N5").Value = 3.451417
$ws.Range("O5").Value = 0.09284703398713583
$ws.Range("P5").Value = 0.09284703398713583
$ws.Range("Q5").Value = 0.1573715765135555
$ws.Range("R5").Value = 1.416344188622
$ws.Range("S5").Value = 0.07262931486557347
$ws.Range("T5").Value = 0.07262931486557347

# Row 6
$ws.Range("G6").Value = 0.03807766666666667
$ws.Range("H6").Value = 0.114233
$ws.Range("I6").Value = 0.2177529884731004
$ws.Range("J6").Value = 0.2177529884731004
$ws.Range("M6").Value = 3.704480666666667
$ws.Range("N6").Value = 11.113442
$ws.Range("O6").Value = 0.298964201395561
$ws.Range("P6").Value = 0.2989642013955609
$ws.Range("Q6").Value = 0.1410579799984444
$ws.Range("R6").Value = 1.269521819986
$ws.Range("S6").Value = 0.06510034830035726
$ws.Range("T6").Value = 0.06510034830035724

# Row 7
$ws.Range("G7").Value = 0.03807766666666667
$ws.Range("H7").Value = 0.114233
$ws.Range("I7").Value = 0.2177529884731004
$ws.Range("J7").Value = 0.2177529884731004
$ws.Range("O7").Value = 0.3765624616238499
$ws.Range("P7").Value = 0.3765624616238499
$ws.Range("Q7").Value = 0.1776705703624445
$ws.Range("R7").Value = 1.599035133262
$ws.Range("S7").Value = 0.08199760136538051
$ws.Range("T7").Value = 0.08199760136538051

# Row 8
$ws.Range("G8").Value = 0.03807766666666667
$ws.Range("H8").Value = 0.114233
$ws.Range("I8").Value = 0.2177529884731004
$ws.Range("J8").Value = 0.2177529884731004
$ws.Range("M8").Value = 2.870093333333334
$ws.Range("N8").Value = 8.610280000000001
$ws.Range("O8").Value = 0.2316263029934534
$ws.Range("P8").Value = 0.2316263029934534
$ws.Range("Q8").Value = 0.1092864572488889
$ws.Range("R8").Value = 0.9835781152400002
$ws.Range("S8").Value = 0.05043731968580032
$ws.Range("T8").Value = 0.05043731968580032

# Row 9
$ws.Range("G9").Value = 0.03807766666666667
$ws.Range("H9").Value = 0.114233
$ws.Range("I9").Value = 0.2177529884731004
$ws.Range("J9").Value = 0.2177529884731004
$ws.Range("M9").Value = 1.150472333333333
$ws.Range("N9").Value = 3.451417
$ws.Range("O9").Value = 0.09284703398713583
$ws.Range("P9").Value = 0.09284703398713583
$ws.Range("Q9").Value = 0.04380730201788889
$ws.Range("R9").Value = 0.394265718161
$ws.Range("S9").Value = 0.02021771912156235
$ws.Range("T9").Value = 0.02021771912156235

